# Mifos -> Finflux: insert a new (blank) column before column N on the
# "Repayment schedule" sheet, shifting the old N/O/P ("Late" / blank /
# "Outstanding") columns one place to the right (-> O/P/Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

$ws.Columns("N").Insert() | Out-Null

$ws.Range("R6").Select() | Out-Null
